# Generate Report for Handback
# Updates the localization-status report after the handback for
# 65f09445-3d6a-40c8-8e93-e3c21a26c086 completes: status flips from
# "Ready for handoff" to "Handed back: in sync with en-US", the
# "Latest Handback DateTime" is refreshed per-locale, and the stale
# "out of date" error detail is cleared now that the handback is current.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for 65f09445-...md (row 3) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row for 65f09445-...md (row 3) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-09-06 09:05:57"
$zhcn.Range("P3").Value = ""
$zhcn.Columns.Item(16).ColumnWidth = 12.86

# --- de-de sheet: row for 65f09445-...md (row 3) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-09-06 09:06:20"
$dede.Range("P3").Value = ""
$dede.Columns.Item(16).ColumnWidth = 12.86
